$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Linear Regression row (row 2) metrics
$ws.Range("B2").Value = 0.9862
$ws.Range("C2").Value = 0.6697
$ws.Range("D2").Value = 0.5242

# Add new row 3: PolynomialRegression
$ws.Range("A3").Value = "PolynomialRegression"
$ws.Range("B3").Value = 0.9867
$ws.Range("C3").Value = 0.6233
$ws.Range("D3").Value = 0.5562

# Add new row 4: Random Forest
$ws.Range("A4").Value = "Random Forest"
$ws.Range("B4").Value = 0.9584
$ws.Range("C4").Value = 1.8186
$ws.Range("D4").Value = 1.0068
